$d = $word.ActiveDocument
$xml = $d.WordOpenXML

# 1) Title style: drop direct spacing/kern from rPr
$old1 = '<w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:spacing w:val="-10"/><w:kern w:val="28"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr>'
$new1 = '<w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:eastAsiaTheme="majorEastAsia" w:hAnsiTheme="majorHAnsi" w:cstheme="majorBidi"/><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr>'
$count1 = ([regex]::Matches($xml, [regex]::Escape($old1))).Count
$xml = $xml.Replace($old1, $new1)

# 2) Author style: add basedOn Title, drop jc, add rPr sz/szCs 24
$old2 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Author"><w:name w:val="Author"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/><w:jc w:val="center"/></w:pPr></w:style>'
$new2 = '<w:style w:type="paragraph" w:customStyle="1" w:styleId="Author"><w:name w:val="Author"/><w:basedOn w:val="Title"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/></w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:style>'
$xml = $xml.Replace($old2, $new2)

# 3) Date style: add basedOn Title, drop jc, add rPr sz/szCs 24
$old3 = '<w:style w:type="paragraph" w:styleId="Date"><w:name w:val="Date"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/><w:jc w:val="center"/></w:pPr></w:style>'
$new3 = '<w:style w:type="paragraph" w:styleId="Date"><w:name w:val="Date"/><w:basedOn w:val="Title"/><w:next w:val="BodyText"/><w:qFormat/><w:pPr><w:keepNext/><w:keepLines/></w:pPr><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:style>'
$xml = $xml.Replace($old3, $new3)

$d.WordOpenXML = $xml
Write-Host ("Title/TitleChar replacements: {0}" -f $count1)
